# Part 2 Step 3: Updated the test plan to have the correct values for the new tests
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Rows 7-10: update Method Inputs (F) / Expected Result (G) text for the
# existing __init__ test cases now that is_borrowed/item_id are part of the
# constructor signature.
# ---------------------------------------------------------------------------
$ws.Range("F7").Value = "title = `"Book Name`"                                                author = `"Author Name`"                                       genre = FICTION                                                               item_id = 1                                                                               is_borrowed = true"
$ws.Range("G7").Value = "The instnace is inititalized correctly, no errors."
$ws.Range("F8").Value = "title = `"`"                                                                                     author = `"Author Name`"                                         genre = NON_FICTION                                                       item_id = 1                                                                               is_borrowed = true"
$ws.Range("G8").Value = "ValueError(`"Title cannot be blank.`")"
$ws.Range("F9").Value = "title = `"Book Name`"                                                                                     author = `"`"                                                         genre = TRUE_CRIME                                                             item_id = 1                                                                               is_borrowed = true"
$ws.Range("G9").Value = "raise ValueError(`"Author cannot be blank.`")"
$ws.Range("F10").Value = "title = `"Book Name`"                                                                                     author = `"Author Name`"                                         genre = `"This is wrong`"                                             item_id = 1                                                                               is_borrowed = true"
$ws.Range("G10").Value = "ValueError(`"Invalid Genre`")"

# ---------------------------------------------------------------------------
# Rows 11-13: the attribute-getter tests (title/author/Genre) now also carry
# a Preconditions value and a "None" Method Inputs value. G13's expected
# result changes from FICTION to genre.FICTION.
# ---------------------------------------------------------------------------
$ws.Range("E11").Value = "The object is initialized correctly     title = `"Book Name`"                                                author = `"Author Name`"                                       genre = FICTION"
$ws.Range("F11").Value = "None"
$ws.Range("G11").Value = "`"Book Name`""
$ws.Range("E12").Value = "The object is initialized correctly     title = `"Book Name`"                                                author = `"Author Name`"                                       genre = FICTION"
$ws.Range("F12").Value = "None"
$ws.Range("G12").Value = "`"Author Name`""
$ws.Range("E13").Value = "The object is initialized correctly     title = `"Book Name`"                                                author = `"Author Name`"                                       genre = FICTION"
$ws.Range("F13").Value = "None"
$ws.Range("G13").Value = "genre.FICTION"

# ---------------------------------------------------------------------------
# Rows 14-17: brand new test cases for item_id / is_barrowed getters plus the
# __init__ validation tests for non-int item_id and non-bool is_borrowed.
# ---------------------------------------------------------------------------
$ws.Range("B14").Value = 8
$ws.Range("C14").Value = "item_id"
$ws.Range("D14").Value = "returns the item id"
$ws.Range("E14").Value = "The object is initialized correctly     title = `"Book Name`"                                                author = `"Author Name`"                                       genre = FICTION                                                  item_id = 1                                                               is_borrowed = true"
$ws.Range("F14").Value = "None"
$ws.Range("G14").Value = 1
$ws.Range("G14").HorizontalAlignment = -4131
$ws.Range("G14").VerticalAlignment = -4160
$ws.Range("G14").WrapText = $true

$ws.Range("B15").Value = 9
$ws.Range("C15").Value = "is_barrowed"
$ws.Range("D15").Value = "returns true or false "
$ws.Range("E15").Value = "The object is initialized correctly     title = `"Book Name`"                                                author = `"Author Name`"                                       genre = FICTION                                                  item_id = 1                                                               is_borrowed = true"
$ws.Range("F15").Value = "None"
$ws.Range("G15").Value = $true
$ws.Range("G15").HorizontalAlignment = -4131
$ws.Range("G15").VerticalAlignment = -4160
$ws.Range("G15").WrapText = $true

$ws.Range("B16").Value = 10
$ws.Range("C16").Value = "__init__"
$ws.Range("D16").Value = "exception raised when item_id is not a int"
$ws.Range("E16").Value = "The object is initialized correctly     title = `"Book Name`"                                                author = `"Author Name`"                                       genre = FICTION                                                  item_id = 1                                                               is_borrowed = true"
$ws.Range("F16").Value = "title = `"Book Name`"                                                author = `"Author Name`"                                       genre = FICTION                                                               item_id = `"Invalid input`"                                                                               is_borrowed = true"
$ws.Range("G16").Value = "ValueError(`"Item Id must be numeric.`")"

$ws.Range("B17").Value = 11
$ws.Range("C17").Value = "__init__"
$ws.Range("D17").Value = "exception raised when is_barrowed is not a bool"
$ws.Range("E17").Value = "The object is initialized correctly     title = `"Book Name`"                                                author = `"Author Name`"                                       genre = FICTION                                                  item_id = 1                                                               is_borrowed = true"
$ws.Range("F17").Value = "title = `"Book Name`"                                                author = `"Author Name`"                                       genre = FICTION                                                               item_id = 1                                                                               is_borrowed = `"Invalid input`""
$ws.Range("G17").Value = "ValueError(`"Is Borrowed must be a boolean value.`")"

# ---------------------------------------------------------------------------
# Rows 18-24: the blank placeholder rows shift up by two Test Case IDs since
# two new real test cases (16 and 17) were inserted above them.
# ---------------------------------------------------------------------------
$ws.Range("B18").Value = 12
$ws.Range("B19").Value = 13
$ws.Range("B20").Value = 14
$ws.Range("B21").Value = 15
$ws.Range("B22").Value = 16
$ws.Range("B23").Value = 17
$ws.Range("B24").Value = 18

# ---------------------------------------------------------------------------
# Row heights: rows 7/8 now match the shorter 72pt rows, and the new data
# rows 14-17 use the 86.4pt height shared by the other populated rows.
# ---------------------------------------------------------------------------
$ws.Rows(7).RowHeight = 72
$ws.Rows(8).RowHeight = 72
$ws.Rows(14).RowHeight = 86.4
$ws.Rows(15).RowHeight = 86.4
$ws.Rows(16).RowHeight = 86.4
$ws.Rows(17).RowHeight = 86.4

# ---------------------------------------------------------------------------
# View state: selection moved to G16 with the window scrolled down.
# ---------------------------------------------------------------------------
$excel.ActiveWindow.ScrollRow = 14
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("G16").Select()

# Page orientation is explicitly portrait.
$ws.PageSetup.Orientation = 1
